$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, shifting rows 11:16 down to 12:17.
$ws.Rows.Item(11).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Copy the style/format from the (now shifted) row 12 -- which is the original row 11 --
# up into the freshly inserted row 11, so the new row inherits the same look.
$ws.Range("A12:Q12").Copy()
$ws.Range("A11:Q11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "0"
$ws.Range("B11").Value = "villageScene1"
$ws.Range("C11").Value = "villageScene1"
$ws.Range("D11").Value = 100
$ws.Range("E11").Value = 500000
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "89,104,0"
$ws.Range("H11").Value = 500
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "Sources/Music/Town"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "Sources/Music/Town"
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").NumberFormat = "@"
$ws.Range("N11").Value = "UI/ChronoBlade_forest_wallpaper"
$ws.Range("O11").NumberFormat = "@"
$ws.Range("O11").Value = "0,8,7"
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "45,180"
$ws.Range("Q11").NumberFormat = "@"
$ws.Range("Q11").Value = "../NFDataCfg/Ini/Navigation/srv_demo.navmesh"

$ws.Range("L10").Select()
